$d = $word.ActiveDocument

# Locate the paragraph that contains the metric we need to update:
# "Dia 12/09: 2hr e 30min (1 dia)"  ->  "Dia 12/09: 3hr (1 dia)"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.Contains("2hr e 30min")) {
        $target = $cand
    }
}

$p = $target
$pStart = $p.Range.Start
$full = $p.Range.Text

# Replace "2hr e 30min" with "3hr". This naturally collapses the runs that
# used to spell out "2" / "hr" / " e " / "30" / "min" into a single run.
$i2hr = $full.IndexOf("2hr e 30min")
$metricStart = $pStart + $i2hr
$metricEnd = $metricStart + "2hr e 30min".Length
$metricRange = $d.Range($metricStart, $metricEnd)
$metricRange.Text = "3hr"

# Re-derive the paragraph boundaries/text after the replacement above, then
# split "Dia 12/09:" away from " 3hr" into their own runs (matching the
# target formatting), while leaving the trailing " (1 dia)" run untouched.
$full2 = $p.Range.Text
$pStart2 = $p.Range.Start

$iDia = $full2.IndexOf("Dia 12/09:")
$colonEnd = $pStart2 + $iDia + "Dia 12/09:".Length
$threeHrEnd = $colonEnd + " 3hr".Length

$r1 = $d.Range($pStart2, $colonEnd)
$r2 = $d.Range($colonEnd, $threeHrEnd)

# Forcing (and immediately reverting) a character-formatting property on each
# sub-range makes the engine keep them as distinct runs instead of merging
# them back with their identically-formatted neighbours.
$r1.Bold = 1
$r1.Bold = 0

$r2.Bold = 1
$r2.Bold = 0

Write-Output $p.Range.Text
